# Auto-generated Excel COM-interop script
# Applies numeric cell-value corrections to the "Siren_Profits" workbook sheets
# (profit/price recalculations) per the authoritative diff.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Cells.Item(9, 8).Value = 41667124  # H9: was 50000520
$ws.Cells.Item(9, 10).Value = 691.8570999999999  # J9: was 918.6
$ws.Cells.Item(9, 12).Value = 691.8570999999999  # L9: was 918.6
$ws.Cells.Item(9, 14).Value = -1029.8571  # N9: was -1256.6
# Row 55
$ws.Cells.Item(55, 8).Value = 971.8461  # H55: was 1138.3636
$ws.Cells.Item(55, 9).Value = 104  # I55: was 120
$ws.Cells.Item(55, 11).Value = 104  # K55: was 120
$ws.Cells.Item(55, 13).Value = 110  # M55: was 94
# Row 69
$ws.Cells.Item(69, 8).Value = 16497.5  # H69: was 14997
$ws.Cells.Item(69, 9).Value = 14995  # I69: was 14997
$ws.Cells.Item(69, 10).Value = 18000  # J69: was 0
$ws.Cells.Item(69, 11).Value = 44985  # K69: was 44991
$ws.Cells.Item(69, 12).Value = 54000  # L69: was 0
$ws.Cells.Item(69, 13).Value = -44111  # M69: was -44117
$ws.Cells.Item(69, 14).Value = -55748  # N69: was NEWCELL
# Row 72
$ws.Cells.Item(72, 8).Value = 16497.5  # H72: was 14997
$ws.Cells.Item(72, 9).Value = 14995  # I72: was 14997
$ws.Cells.Item(72, 10).Value = 18000  # J72: was 0
$ws.Cells.Item(72, 11).Value = 134955  # K72: was 134973
$ws.Cells.Item(72, 12).Value = 162000  # L72: was 0
$ws.Cells.Item(72, 13).Value = -130587  # M72: was -130605
$ws.Cells.Item(72, 14).Value = -170736  # N72: was NEWCELL
# Row 74
$ws.Cells.Item(74, 8).Value = 4407.2856  # H74: was 4011.2778
$ws.Cells.Item(74, 9).Value = 2740.2  # I74: was 2707.3572
$ws.Cells.Item(74, 11).Value = 2740.2  # K74: was 2707.3572
$ws.Cells.Item(74, 13).Value = -1804.2  # M74: was -1771.3572
# Row 77
$ws.Cells.Item(77, 8).Value = 4407.2856  # H77: was 4011.2778
$ws.Cells.Item(77, 9).Value = 2740.2  # I77: was 2707.3572
$ws.Cells.Item(77, 11).Value = 13701  # K77: was 13536.786
$ws.Cells.Item(77, 13).Value = -9021  # M77: was -8856.786
# Row 95
$ws.Cells.Item(95, 8).Value = 21541.334  # H95: was 29624
$ws.Cells.Item(95, 10).Value = 21541.334  # J95: was 29624
$ws.Cells.Item(95, 12).Value = 21541.334  # L95: was 29624
$ws.Cells.Item(95, 14).Value = -27033.334  # N95: was -35116
# Row 112
$ws.Cells.Item(112, 8).Value = 28171.682  # H112: was 26556.256
$ws.Cells.Item(112, 10).Value = 28718.89  # J112: was 27027.38
$ws.Cells.Item(112, 12).Value = 86156.67  # L112: was 81082.14
$ws.Cells.Item(112, 14).Value = -88372.67  # N112: was -83298.14
# Row 121
$ws.Cells.Item(121, 8).Value = 1081.2858  # H121: was 1086.0476
$ws.Cells.Item(121, 10).Value = 1081.2858  # J121: was 1086.0476
$ws.Cells.Item(121, 12).Value = 3243.8574  # L121: was 3258.142800000001
$ws.Cells.Item(121, 14).Value = -6737.857400000001  # N121: was -6752.142800000001
# Row 129
$ws.Cells.Item(129, 8).Value = 76925030  # H129: was 83335170
$ws.Cells.Item(129, 9).Value = 2164.6667  # I129: was 1863.5
$ws.Cells.Item(129, 10).Value = 142858910  # J129: was 166668460
$ws.Cells.Item(129, 11).Value = 6494.000100000001  # K129: was 5590.5
$ws.Cells.Item(129, 12).Value = 428576730  # L129: was 500005380
$ws.Cells.Item(129, 13).Value = -1494.000100000001  # M129: was -590.5
$ws.Cells.Item(129, 14).Value = -428586730  # N129: was -500015380
# Row 131
$ws.Cells.Item(131, 8).Value = 3192.3076  # H131: was 3257.1428
$ws.Cells.Item(131, 10).Value = 8033.3335  # J131: was 7050
$ws.Cells.Item(131, 12).Value = 24100.0005  # L131: was 21150
$ws.Cells.Item(131, 14).Value = -34180.00049999999  # N131: was -31230
# Row 137
$ws.Cells.Item(137, 8).Value = 801959.3  # H137: was 947624.4
$ws.Cells.Item(137, 9).Value = 1719062.9  # I137: was 2578194.2
$ws.Cells.Item(137, 10).Value = 15870.571  # J137: was 15870.143
$ws.Cells.Item(137, 11).Value = 5157188.699999999  # K137: was 7734582.600000001
$ws.Cells.Item(137, 12).Value = 47611.713  # L137: was 47610.429
$ws.Cells.Item(137, 13).Value = -5154638.699999999  # M137: was -7732032.600000001
$ws.Cells.Item(137, 14).Value = -52711.713  # N137: was -52710.429
# Row 138
$ws.Cells.Item(138, 8).Value = 5298.1606  # H138: was 5391.0386
$ws.Cells.Item(138, 10).Value = 5874.913  # J138: was 6010.894
$ws.Cells.Item(138, 12).Value = 17624.739  # L138: was 18032.682
$ws.Cells.Item(138, 14).Value = -27904.739  # N138: was -28312.682
# Row 139
$ws.Cells.Item(139, 8).Value = 142088.6  # H139: was 141166.17
$ws.Cells.Item(139, 10).Value = 142088.6  # J139: was 141166.17
$ws.Cells.Item(139, 12).Value = 142088.6  # L139: was 141166.17
$ws.Cells.Item(139, 14).Value = -152368.6  # N139: was -151446.17
# Row 141
$ws.Cells.Item(141, 8).Value = 3570.8147  # H141: was 2973.853
$ws.Cells.Item(141, 9).Value = 2043.45  # I141: was 1772.52
$ws.Cells.Item(141, 10).Value = 7934.7144  # J141: was 6310.8887
$ws.Cells.Item(141, 11).Value = 6130.35  # K141: was 5317.559999999999
$ws.Cells.Item(141, 12).Value = 23804.1432  # L141: was 18932.6661
$ws.Cells.Item(141, 13).Value = -950.3500000000004  # M141: was -137.5599999999995
$ws.Cells.Item(141, 14).Value = -34164.1432  # N141: was -29292.6661

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 4170  # H32: was 4587.086
$ws.Cells.Item(32, 9).Value = 4237.6577  # I32: was 4674.9707
$ws.Cells.Item(32, 11).Value = 4237.6577  # K32: was 4674.9707
$ws.Cells.Item(32, 13).Value = -3950.6577  # M32: was -4387.9707
# Row 45
$ws.Cells.Item(45, 8).Value = 110404.85  # H45: was 67975.664
$ws.Cells.Item(45, 10).Value = 4650.6665  # J45: was 3316
$ws.Cells.Item(45, 12).Value = 4650.6665  # L45: was 3316
$ws.Cells.Item(45, 14).Value = -5404.6665  # N45: was -4070
# Row 63
$ws.Cells.Item(63, 8).Value = 2349.75  # H63: was 7054.364
$ws.Cells.Item(63, 9).Value = 2349.75  # I63: was 7054.364
$ws.Cells.Item(63, 11).Value = 2349.75  # K63: was 7054.364
$ws.Cells.Item(63, 13).Value = -1663.75  # M63: was -6368.364
# Row 66
$ws.Cells.Item(66, 8).Value = 2349.75  # H66: was 7054.364
$ws.Cells.Item(66, 9).Value = 2349.75  # I66: was 7054.364
$ws.Cells.Item(66, 11).Value = 11748.75  # K66: was 35271.82
$ws.Cells.Item(66, 13).Value = -8316.75  # M66: was -31839.82
# Row 102
$ws.Cells.Item(102, 8).Value = 5417.0713  # H102: was 3073.3845
$ws.Cells.Item(102, 9).Value = 5237  # I102: was 2788.0417
$ws.Cells.Item(102, 11).Value = 5237  # K102: was 2788.0417
$ws.Cells.Item(102, 13).Value = -3615  # M102: was -1166.0417
# Row 132
$ws.Cells.Item(132, 8).Value = 2875.86  # H132: was 2920.0205
$ws.Cells.Item(132, 9).Value = 2245.0789  # I132: was 2286.5134
$ws.Cells.Item(132, 11).Value = 6735.236699999999  # K132: was 6859.540199999999
$ws.Cells.Item(132, 13).Value = -4205.236699999999  # M132: was -4329.540199999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 82
$ws.Cells.Item(82, 8).Value = 57199.6  # H82: was 64499.75
$ws.Cells.Item(82, 9).Value = 28666.334  # I82: was 29000
$ws.Cells.Item(82, 11).Value = 28666.334  # K82: was 29000
$ws.Cells.Item(82, 13).Value = -28283.334  # M82: was -28617
# Row 85
$ws.Cells.Item(85, 8).Value = 57199.6  # H85: was 64499.75
$ws.Cells.Item(85, 9).Value = 28666.334  # I85: was 29000
$ws.Cells.Item(85, 11).Value = 28666.334  # K85: was 29000
$ws.Cells.Item(85, 13).Value = -27340.334  # M85: was -27674
# Row 86
$ws.Cells.Item(86, 8).Value = 60000  # H86: was 16999.2
$ws.Cells.Item(86, 9).Value = 60000  # I86: was 16999.2
$ws.Cells.Item(86, 11).Value = 60000  # K86: was 16999.2
$ws.Cells.Item(86, 13).Value = -58877  # M86: was -15876.2
# Row 89
$ws.Cells.Item(89, 8).Value = 60000  # H89: was 16999.2
$ws.Cells.Item(89, 9).Value = 60000  # I89: was 16999.2
$ws.Cells.Item(89, 11).Value = 300000  # K89: was 84996
$ws.Cells.Item(89, 13).Value = -294384  # M89: was -79380
# Row 105
$ws.Cells.Item(105, 8).Value = 56759.285  # H105: was 54360.953
$ws.Cells.Item(105, 10).Value = 4810  # J105: was 4728.6
$ws.Cells.Item(105, 12).Value = 4810  # L105: was 4728.6
$ws.Cells.Item(105, 14).Value = -8304  # N105: was -8222.6
# Row 107
$ws.Cells.Item(107, 8).Value = 3696.4614  # H107: was 3704.1155
$ws.Cells.Item(107, 9).Value = 3401.1333  # I107: was 3414.4
$ws.Cells.Item(107, 11).Value = 3401.1333  # K107: was 3414.4
$ws.Cells.Item(107, 13).Value = -1481.1333  # M107: was -1494.4

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Cells.Item(7, 8).Value = 177  # H7: was 185.66667
$ws.Cells.Item(7, 9).Value = 177  # I7: was 174
$ws.Cells.Item(7, 10).Value = 0  # J7: was 244
$ws.Cells.Item(7, 11).Value = 177  # K7: was 174
$ws.Cells.Item(7, 12).Value = 0  # L7: was 244
$ws.Cells.Item(7, 13).Value = -64  # M7: was -61
$ws.Cells.Item(7, 14).ClearContents()  # N7: was -470
# Row 31
$ws.Cells.Item(31, 8).Value = 2738.4167  # H31: was 2738.25
$ws.Cells.Item(31, 9).Value = 1702.8572  # I31: was 1820
$ws.Cells.Item(31, 10).Value = 2915.2195  # J31: was 2869.4285
$ws.Cells.Item(31, 11).Value = 1702.8572  # K31: was 1820
$ws.Cells.Item(31, 12).Value = 2915.2195  # L31: was 2869.4285
$ws.Cells.Item(31, 13).Value = -1407.8572  # M31: was -1525
$ws.Cells.Item(31, 14).Value = -3505.2195  # N31: was -3459.4285
# Row 34
$ws.Cells.Item(34, 8).Value = 2738.4167  # H34: was 2738.25
$ws.Cells.Item(34, 9).Value = 1702.8572  # I34: was 1820
$ws.Cells.Item(34, 10).Value = 2915.2195  # J34: was 2869.4285
$ws.Cells.Item(34, 11).Value = 1702.8572  # K34: was 1820
$ws.Cells.Item(34, 12).Value = 2915.2195  # L34: was 2869.4285
$ws.Cells.Item(34, 13).Value = -1500.8572  # M34: was -1618
$ws.Cells.Item(34, 14).Value = -3319.2195  # N34: was -3273.4285
# Row 99
$ws.Cells.Item(99, 8).Value = 13892154  # H99: was 15628199
$ws.Cells.Item(99, 10).Value = 4399.25  # J99: was 4599.6665
$ws.Cells.Item(99, 12).Value = 4399.25  # L99: was 4599.6665
$ws.Cells.Item(99, 14).Value = -7395.25  # N99: was -7595.6665
# Row 104
$ws.Cells.Item(104, 8).Value = 49997.25  # H104: was 53990.8
$ws.Cells.Item(104, 10).Value = 49997.25  # J104: was 53990.8
$ws.Cells.Item(104, 12).Value = 49997.25  # L104: was 53990.8
$ws.Cells.Item(104, 14).Value = -55239.25  # N104: was -59232.8
# Row 107
$ws.Cells.Item(107, 8).Value = 11385.903  # H107: was 11061.344
$ws.Cells.Item(107, 9).Value = 21104.312  # I107: was 19921.705
$ws.Cells.Item(107, 11).Value = 21104.312  # K107: was 19921.705
$ws.Cells.Item(107, 13).Value = -19184.312  # M107: was -18001.705
# Row 126
$ws.Cells.Item(126, 8).Value = 13892154  # H126: was 15628199
$ws.Cells.Item(126, 10).Value = 4399.25  # J126: was 4599.6665
$ws.Cells.Item(126, 12).Value = 13197.75  # L126: was 13798.9995
$ws.Cells.Item(126, 14).Value = -18137.75  # N126: was -18738.9995
# Row 134
$ws.Cells.Item(134, 8).Value = 2506462.5  # H134: was 2506523.5
$ws.Cells.Item(134, 10).Value = 2992.25  # J134: was 3182.25
$ws.Cells.Item(134, 12).Value = 8976.75  # L134: was 9546.75
$ws.Cells.Item(134, 14).Value = -14046.75  # N134: was -14616.75

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Cells.Item(2, 8).Value = 80.5  # H2: was 80.81579000000001
$ws.Cells.Item(2, 9).Value = 92.454544  # I2: was 90.21738999999999
$ws.Cells.Item(2, 10).Value = 64.0625  # J2: was 66.40000000000001
$ws.Cells.Item(2, 11).Value = 554.727264  # K2: was 541.3043399999999
$ws.Cells.Item(2, 12).Value = 384.375  # L2: was 398.4
$ws.Cells.Item(2, 13).Value = -441.727264  # M2: was -428.3043399999999
$ws.Cells.Item(2, 14).Value = -610.375  # N2: was -624.4000000000001
# Row 33
$ws.Cells.Item(33, 8).Value = 228.75  # H33: was 206.25
$ws.Cells.Item(33, 9).Value = 190  # I33: was 164.28572
$ws.Cells.Item(33, 11).Value = 1140  # K33: was 985.71432
$ws.Cells.Item(33, 13).Value = -857  # M33: was -702.71432
# Row 124
$ws.Cells.Item(124, 8).Value = 4238.222  # H124: was 5160.1113
$ws.Cells.Item(124, 9).Value = 1047.25  # I124: was 1264.3334
$ws.Cells.Item(124, 10).Value = 6791  # J124: was 7108
$ws.Cells.Item(124, 11).Value = 3141.75  # K124: was 3793.0002
$ws.Cells.Item(124, 12).Value = 20373  # L124: was 21324
$ws.Cells.Item(124, 13).Value = 1768.25  # M124: was 1116.9998
$ws.Cells.Item(124, 14).Value = -30193  # N124: was -31144
# Row 131
$ws.Cells.Item(131, 8).Value = 34486452  # H131: was 35718076
$ws.Cells.Item(131, 9).Value = 50004268  # I131: was 52636024
$ws.Cells.Item(131, 11).Value = 150012804  # K131: was 157908072
$ws.Cells.Item(131, 13).Value = -150007764  # M131: was -157903032
# Row 137
$ws.Cells.Item(137, 8).Value = 10153.267  # H137: was 10519.214
$ws.Cells.Item(137, 9).Value = 5671.5  # I137: was 5799.8
$ws.Cells.Item(137, 11).Value = 17014.5  # K137: was 17399.4
$ws.Cells.Item(137, 13).Value = -11914.5  # M137: was -12299.4
# Row 140
$ws.Cells.Item(140, 8).Value = 2472.75  # H140: was 2675.3076
$ws.Cells.Item(140, 9).Value = 2472.75  # I140: was 2675.3076
$ws.Cells.Item(140, 11).Value = 7418.25  # K140: was 8025.9228
$ws.Cells.Item(140, 13).Value = -2238.25  # M140: was -2845.9228

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 105
$ws.Cells.Item(105, 8).Value = 55999.25  # H105: was 61332.332
$ws.Cells.Item(105, 10).Value = 55999.25  # J105: was 61332.332
$ws.Cells.Item(105, 12).Value = 55999.25  # L105: was 61332.332
$ws.Cells.Item(105, 14).Value = -62987.25  # N105: was -68320.33199999999
# Row 126
$ws.Cells.Item(126, 8).Value = 13122.808  # H126: was 13801.167
$ws.Cells.Item(126, 9).Value = 12452.294  # I126: was 13024.375
$ws.Cells.Item(126, 10).Value = 14389.333  # J126: was 15354.75
$ws.Cells.Item(126, 11).Value = 37356.882  # K126: was 39073.125
$ws.Cells.Item(126, 12).Value = 43167.999  # L126: was 46064.25
$ws.Cells.Item(126, 13).Value = -34886.882  # M126: was -36603.125
$ws.Cells.Item(126, 14).Value = -48107.999  # N126: was -51004.25
# Row 132
$ws.Cells.Item(132, 8).Value = 3634.1667  # H132: was 3881.6365
$ws.Cells.Item(132, 9).Value = 1861.5  # I132: was 1967
$ws.Cells.Item(132, 11).Value = 5584.5  # K132: was 5901
$ws.Cells.Item(132, 13).Value = -3054.5  # M132: was -3371

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Cells.Item(16, 8).Value = 1797.9231  # H16: was 1804.0714
$ws.Cells.Item(16, 9).Value = 1781.0834  # I16: was 1789
$ws.Cells.Item(16, 11).Value = 1781.0834  # K16: was 1789
$ws.Cells.Item(16, 13).Value = -1611.0834  # M16: was -1619
# Row 33
$ws.Cells.Item(33, 8).Value = 18998  # H33: was 3766872
$ws.Cells.Item(33, 9).Value = 12499.5  # I33: was 7511248.5
$ws.Cells.Item(33, 10).Value = 23330.334  # J33: was 22495.75
$ws.Cells.Item(33, 11).Value = 12499.5  # K33: was 7511248.5
$ws.Cells.Item(33, 12).Value = 23330.334  # L33: was 22495.75
$ws.Cells.Item(33, 13).Value = -12209.5  # M33: was -7510958.5
$ws.Cells.Item(33, 14).Value = -23910.334  # N33: was -23075.75
# Row 38
$ws.Cells.Item(38, 8).Value = 40504.168  # H38: was 39574.57
$ws.Cells.Item(38, 10).Value = 45799  # J38: was 43832
$ws.Cells.Item(38, 12).Value = 45799  # L38: was 43832
$ws.Cells.Item(38, 14).Value = -46619  # N38: was -44652
# Row 44
$ws.Cells.Item(44, 8).Value = 15996.667  # H44: was 15999.333
$ws.Cells.Item(44, 10).Value = 15996.667  # J44: was 15999.333
$ws.Cells.Item(44, 12).Value = 15996.667  # L44: was 15999.333
$ws.Cells.Item(44, 14).Value = -16908.667  # N44: was -16911.333
# Row 50
$ws.Cells.Item(50, 8).Value = 48357  # H50: was 66662.336
$ws.Cells.Item(50, 9).Value = 25076  # I50: was 60000
$ws.Cells.Item(50, 10).Value = 59997.5  # J50: was 69993.5
$ws.Cells.Item(50, 11).Value = 25076  # K50: was 60000
$ws.Cells.Item(50, 12).Value = 59997.5  # L50: was 69993.5
$ws.Cells.Item(50, 13).Value = -24439  # M50: was -59363
$ws.Cells.Item(50, 14).Value = -61271.5  # N50: was -71267.5
# Row 53
$ws.Cells.Item(53, 8).Value = 9499.5  # H53: was 12665.333
$ws.Cells.Item(53, 10).Value = 6000  # J53: was 12498.5
$ws.Cells.Item(53, 12).Value = 6000  # L53: was 12498.5
$ws.Cells.Item(53, 14).Value = -7036  # N53: was -13534.5
# Row 54
$ws.Cells.Item(54, 8).Value = 0  # H54: was 25042
$ws.Cells.Item(54, 10).Value = 0  # J54: was 25042
$ws.Cells.Item(54, 12).Value = 0  # L54: was 25042
$ws.Cells.Item(54, 14).ClearContents()  # N54: was -26330
# Row 56
$ws.Cells.Item(56, 8).Value = 15025.5  # H56: was 16349.333
$ws.Cells.Item(56, 10).Value = 0  # J56: was 18997
$ws.Cells.Item(56, 12).Value = 0  # L56: was 18997
$ws.Cells.Item(56, 14).ClearContents()  # N56: was -20379
# Row 60
$ws.Cells.Item(60, 8).Value = 23250  # H60: was 44332.332
$ws.Cells.Item(60, 10).Value = 23250  # J60: was 44332.332
$ws.Cells.Item(60, 12).Value = 23250  # L60: was 44332.332
$ws.Cells.Item(60, 14).Value = -24268  # N60: was -45350.332
# Row 61
$ws.Cells.Item(61, 8).Value = 3193.1538  # H61: was 3333.4167
$ws.Cells.Item(61, 9).Value = 3042.5833  # I61: was 3251.1
$ws.Cells.Item(61, 10).Value = 5000  # J61: was 3745
$ws.Cells.Item(61, 11).Value = 3042.5833  # K61: was 3251.1
$ws.Cells.Item(61, 12).Value = 5000  # L61: was 3745
$ws.Cells.Item(61, 13).Value = -2840.5833  # M61: was -3049.1
$ws.Cells.Item(61, 14).Value = -5404  # N61: was -4149
# Row 93
$ws.Cells.Item(93, 8).Value = 6098.091  # H93: was 6547.7
$ws.Cells.Item(93, 9).Value = 10079  # I93: was 12198.25
$ws.Cells.Item(93, 11).Value = 10079  # K93: was 12198.25
$ws.Cells.Item(93, 13).Value = -8831  # M93: was -10950.25
# Row 113
$ws.Cells.Item(113, 8).Value = 3193.1538  # H113: was 3333.4167
$ws.Cells.Item(113, 9).Value = 3042.5833  # I113: was 3251.1
$ws.Cells.Item(113, 10).Value = 5000  # J113: was 3745
$ws.Cells.Item(113, 11).Value = 3042.5833  # K113: was 3251.1
$ws.Cells.Item(113, 12).Value = 5000  # L113: was 3745
$ws.Cells.Item(113, 13).Value = -872.5832999999998  # M113: was -1081.1
$ws.Cells.Item(113, 14).Value = -9340  # N113: was -8085
# Row 122
$ws.Cells.Item(122, 8).Value = 4931.8667  # H122: was 4255.923
$ws.Cells.Item(122, 9).Value = 4007.9  # I122: was 3785.5833
$ws.Cells.Item(122, 10).Value = 6779.8  # J122: was 9900
$ws.Cells.Item(122, 11).Value = 12023.7  # K122: was 11356.7499
$ws.Cells.Item(122, 12).Value = 20339.4  # L122: was 29700
$ws.Cells.Item(122, 13).Value = -9573.700000000001  # M122: was -8906.749899999999
$ws.Cells.Item(122, 14).Value = -25239.4  # N122: was -34600

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 45
$ws.Cells.Item(45, 8).Value = 13710.546  # H45: was 13570.615
$ws.Cells.Item(45, 10).Value = 13916.667  # J45: was 13713.818
$ws.Cells.Item(45, 12).Value = 13916.667  # L45: was 13713.818
$ws.Cells.Item(45, 14).Value = -14898.667  # N45: was -14695.818
# Row 113
$ws.Cells.Item(113, 8).Value = 3566.3572  # H113: was 3352.6
$ws.Cells.Item(113, 9).Value = 2382  # I113: was 2189.4285
$ws.Cells.Item(113, 11).Value = 7146  # K113: was 6568.2855
$ws.Cells.Item(113, 13).Value = -4976  # M113: was -4398.2855
# Row 122
$ws.Cells.Item(122, 8).Value = 7008.1875  # H122: was 7270.8
$ws.Cells.Item(122, 9).Value = 3642  # I122: was 3723.8572
$ws.Cells.Item(122, 11).Value = 10926  # K122: was 11171.5716
$ws.Cells.Item(122, 13).Value = -8476  # M122: was -8721.571599999999
